$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting NCTId..intervention_type one column right.
$ws.Range("C1").EntireColumn.Insert()

# Copy the header style from the (now shifted) D1 header cell onto the new C1 header cell.
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header cell.
$ws.Range("C1").Value = "statut_name"

# New data cell for row 2.
$ws.Range("C2").Value = "pas de résultat ni de publication"
